$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 671505.4
$ws.Range("I64").Value = 912507.4
$ws.Range("J64").Value = 8750
$ws.Range("K64").Value = 912507.4
$ws.Range("L64").Value = 8750
$ws.Range("M64").Value = -912259.4
$ws.Range("N64").Value = -9246

$ws.Range("H67").Value = 671505.4
$ws.Range("I67").Value = 912507.4
$ws.Range("J67").Value = 8750
$ws.Range("K67").Value = 912507.4
$ws.Range("L67").Value = 8750
$ws.Range("M67").Value = -911649.4
$ws.Range("N67").Value = -10466

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19333.967
$ws.Range("I32").Value = 3399.2942
$ws.Range("J32").Value = 93212.91
$ws.Range("K32").Value = 3399.2942
$ws.Range("L32").Value = 93212.91
$ws.Range("M32").Value = -3112.2942
$ws.Range("N32").Value = -93786.91

$ws.Range("H63").Value = 6319.815
$ws.Range("I63").Value = 6810.391
$ws.Range("K63").Value = 6810.391
$ws.Range("M63").Value = -6124.391

$ws.Range("H66").Value = 6319.815
$ws.Range("I66").Value = 6810.391
$ws.Range("K66").Value = 34051.955
$ws.Range("M66").Value = -30619.955

$ws.Range("H74").Value = 3415.5095
$ws.Range("I74").Value = 1062.7632
$ws.Range("J74").Value = 9375.799999999999
$ws.Range("K74").Value = 1062.7632
$ws.Range("L74").Value = 9375.799999999999
$ws.Range("M74").Value = -188.7632000000001
$ws.Range("N74").Value = -11123.8

$ws.Range("H77").Value = 3415.5095
$ws.Range("I77").Value = 1062.7632
$ws.Range("J77").Value = 9375.799999999999
$ws.Range("K77").Value = 5313.816000000001
$ws.Range("L77").Value = 46879
$ws.Range("M77").Value = -945.8160000000007
$ws.Range("N77").Value = -55615

$ws.Range("H124").Value = 31330.6
$ws.Range("J124").Value = 31330.6
$ws.Range("L124").Value = 31330.6
$ws.Range("N124").Value = -41150.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 29000.5
$ws.Range("J35").Value = 8001
$ws.Range("L35").Value = 8001
$ws.Range("N35").Value = -8621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20586

$ws.Range("H31").Value = 1163.9454
$ws.Range("I31").Value = 875.8293
$ws.Range("J31").Value = 2007.7142
$ws.Range("K31").Value = 875.8293
$ws.Range("L31").Value = 2007.7142
$ws.Range("M31").Value = -580.8293
$ws.Range("N31").Value = -2597.7142

$ws.Range("H34").Value = 1163.9454
$ws.Range("I34").Value = 875.8293
$ws.Range("J34").Value = 2007.7142
$ws.Range("K34").Value = 875.8293
$ws.Range("L34").Value = 2007.7142
$ws.Range("M34").Value = -673.8293
$ws.Range("N34").Value = -2411.7142

$ws.Range("H58").Value = 1720.1666
$ws.Range("I58").Value = 894.3333
$ws.Range("J58").Value = 2546
$ws.Range("K58").Value = 894.3333
$ws.Range("L58").Value = 2546
$ws.Range("M58").Value = -691.3333
$ws.Range("N58").Value = -2952

$ws.Range("H107").Value = 287.18182
$ws.Range("J107").Value = 316.375
$ws.Range("L107").Value = 316.375
$ws.Range("N107").Value = -4156.375

$ws.Range("H132").Value = 2171.4849
$ws.Range("I132").Value = 1448.84
$ws.Range("J132").Value = 4429.75
$ws.Range("K132").Value = 4346.52
$ws.Range("L132").Value = 13289.25
$ws.Range("M132").Value = -1816.52
$ws.Range("N132").Value = -18349.25

$ws.Range("H135").Value = 42652.855
$ws.Range("J135").Value = 42652.855
$ws.Range("L135").Value = 42652.855
$ws.Range("N135").Value = -52792.855

$ws.Range("H136").Value = 1720.1666
$ws.Range("I136").Value = 894.3333
$ws.Range("J136").Value = 2546
$ws.Range("K136").Value = 2682.9999
$ws.Range("L136").Value = 7638
$ws.Range("M136").Value = -132.9998999999998
$ws.Range("N136").Value = -12738

$ws.Range("H137").Value = 57500
$ws.Range("J137").Value = 57500
$ws.Range("L137").Value = 57500
$ws.Range("N137").Value = -67700

$ws.Range("H138").Value = 54440
$ws.Range("J138").Value = 54440
$ws.Range("L138").Value = 54440
$ws.Range("N138").Value = -64720

$ws.Range("H140").Value = 68000
$ws.Range("J140").Value = 68000
$ws.Range("L140").Value = 68000
$ws.Range("N140").Value = -78360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1050
$ws.Range("J80").Value = 1100
$ws.Range("L80").Value = 3300
$ws.Range("N80").Value = -5172

$ws.Range("H83").Value = 1050
$ws.Range("J83").Value = 1100
$ws.Range("L83").Value = 9900
$ws.Range("N83").Value = -19260

$ws.Range("H131").Value = 1424.2878
$ws.Range("I131").Value = 463.75
$ws.Range("J131").Value = 1556.7759
$ws.Range("K131").Value = 1391.25
$ws.Range("L131").Value = 4670.3277
$ws.Range("M131").Value = 3648.75
$ws.Range("N131").Value = -14750.3277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2539.6377
$ws.Range("I132").Value = 2234.111
$ws.Range("J132").Value = 3639.5334
$ws.Range("K132").Value = 6702.333
$ws.Range("L132").Value = 10918.6002
$ws.Range("M132").Value = -4172.333
$ws.Range("N132").Value = -15978.6002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9196.154
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 10668.182
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 10668.182
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -11258.182

$ws.Range("H27").Value = 9196.154
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 10668.182
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 10668.182
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -10882.182

$ws.Range("H110").Value = 19667.334
$ws.Range("J110").Value = 19667.334
$ws.Range("L110").Value = 19667.334
$ws.Range("N110").Value = -27847.334

$ws.Range("H122").Value = 3306.6155
$ws.Range("I122").Value = 2658.8
$ws.Range("J122").Value = 3711.5
$ws.Range("K122").Value = 7976.400000000001
$ws.Range("L122").Value = 11134.5
$ws.Range("M122").Value = -5526.400000000001
$ws.Range("N122").Value = -16034.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11360
$ws.Range("J45").Value = 11360
$ws.Range("L45").Value = 11360
$ws.Range("N45").Value = -12342

$ws.Range("H119").Value = 21259.6
$ws.Range("J119").Value = 21259.6
$ws.Range("L119").Value = 21259.6
$ws.Range("N119").Value = -30935.6

$ws.Range("H132").Value = 10871958
$ws.Range("I132").Value = 14287596
$ws.Range("J132").Value = 4017.9092
$ws.Range("K132").Value = 42862788
$ws.Range("L132").Value = 12053.7276
$ws.Range("M132").Value = -42860258
$ws.Range("N132").Value = -17113.7276

Write-Host "Applied all changes."
